$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3122.4138
$ws.Range("I86").Value = 1624.6428
$ws.Range("K86").Value = 1624.6428
$ws.Range("M86").Value = -501.6428000000001
$ws.Range("H89").Value = 3122.4138
$ws.Range("I89").Value = 1624.6428
$ws.Range("K89").Value = 8123.214
$ws.Range("M89").Value = -2507.214
$ws.Range("H137").Value = 3237.5
$ws.Range("J137").Value = 5000
$ws.Range("L137").Value = 15000
$ws.Range("N137").Value = -20100

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1184.8695
$ws.Range("I2").Value = 883.8
$ws.Range("J2").Value = 1749.375
$ws.Range("K2").Value = 883.8
$ws.Range("L2").Value = 1749.375
$ws.Range("M2").Value = -770.8
$ws.Range("N2").Value = -1975.375
$ws.Range("H32").Value = 11773558
$ws.Range("I32").Value = 16396890
$ws.Range("K32").Value = 16396890
$ws.Range("M32").Value = -16396603
$ws.Range("H116").Value = 1184.8695
$ws.Range("I116").Value = 883.8
$ws.Range("J116").Value = 1749.375
$ws.Range("K116").Value = 883.8
$ws.Range("L116").Value = 1749.375
$ws.Range("M116").Value = 1410.2
$ws.Range("N116").Value = -6337.375
$ws.Range("H122").Value = 2249.8333
$ws.Range("I122").Value = 1307.3914
$ws.Range("J122").Value = 3917.2307
$ws.Range("K122").Value = 3922.1742
$ws.Range("L122").Value = 11751.6921
$ws.Range("M122").Value = -1472.1742
$ws.Range("N122").Value = -16651.6921

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 42000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 42000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 42000
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -42226
$ws.Range("H3").Value = 1184.8695
$ws.Range("I3").Value = 883.8
$ws.Range("J3").Value = 1749.375
$ws.Range("K3").Value = 883.8
$ws.Range("L3").Value = 1749.375
$ws.Range("M3").Value = -769.8
$ws.Range("N3").Value = -1977.375
$ws.Range("H80").Value = 2597.7778
$ws.Range("I80").Value = 637.6
$ws.Range("K80").Value = 637.6
$ws.Range("M80").Value = 360.4
$ws.Range("H83").Value = 2597.7778
$ws.Range("I83").Value = 637.6
$ws.Range("K83").Value = 3188
$ws.Range("M83").Value = 1804
$ws.Range("H86").Value = 17030.857
$ws.Range("I86").Value = 2767.111
$ws.Range("K86").Value = 2767.111
$ws.Range("M86").Value = -1644.111
$ws.Range("H89").Value = 17030.857
$ws.Range("I89").Value = 2767.111
$ws.Range("K89").Value = 13835.555
$ws.Range("M89").Value = -8219.555
$ws.Range("H105").Value = 1364
$ws.Range("I105").Value = 1205
$ws.Range("K105").Value = 1205
$ws.Range("M105").Value = 542

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 7724.5
$ws.Range("I105").Value = 2043.2222
$ws.Range("J105").Value = 13405.777
$ws.Range("K105").Value = 2043.2222
$ws.Range("L105").Value = 13405.777
$ws.Range("M105").Value = -296.2221999999999
$ws.Range("N105").Value = -16899.777
$ws.Range("H107").Value = 931.8461
$ws.Range("I107").Value = 535.1905
$ws.Range("K107").Value = 535.1905
$ws.Range("M107").Value = 1384.8095
$ws.Range("H124").Value = 190056.5
$ws.Range("J124").Value = 190056.5
$ws.Range("L124").Value = 190056.5
$ws.Range("N124").Value = -194966.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 101189.5
$ws.Range("I32").Value = 333566.66
$ws.Range("K32").Value = 1000699.98
$ws.Range("M32").Value = -1000416.98
$ws.Range("H98").Value = 2051.1667
$ws.Range("J98").Value = 2218.5
$ws.Range("L98").Value = 6655.5
$ws.Range("N98").Value = -9651.5
$ws.Range("H140").Value = 1421.091
$ws.Range("I140").Value = 1289.9
$ws.Range("J140").Value = 2733
$ws.Range("K140").Value = 3869.7
$ws.Range("L140").Value = 8199
$ws.Range("M140").Value = 1310.3
$ws.Range("N140").Value = -18559

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3251.5862
$ws.Range("I113").Value = 2573.1428
$ws.Range("K113").Value = 2573.1428
$ws.Range("M113").Value = -403.1428000000001
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H122").Value = 2409.6667
$ws.Range("I122").Value = 2245.9524
$ws.Range("K122").Value = 6737.8572
$ws.Range("M122").Value = -4287.8572
$ws.Range("H132").Value = 1401.8918
$ws.Range("I132").Value = 1209.7587
$ws.Range("J132").Value = 2098.375
$ws.Range("K132").Value = 3629.2761
$ws.Range("L132").Value = 6295.125
$ws.Range("M132").Value = -1099.2761
$ws.Range("N132").Value = -11355.125
$ws.Range("H136").Value = 63519.8
$ws.Range("J136").Value = 63519.8
$ws.Range("L136").Value = 190559.4
$ws.Range("N136").Value = -195659.4

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4289.6055
$ws.Range("I7").Value = 3816.0588
$ws.Range("J7").Value = 4672.952
$ws.Range("K7").Value = 3816.0588
$ws.Range("L7").Value = 4672.952
$ws.Range("M7").Value = -3704.0588
$ws.Range("N7").Value = -4896.952
$ws.Range("H122").Value = 4890.24
$ws.Range("I122").Value = 4453.6875
$ws.Range("K122").Value = 13361.0625
$ws.Range("M122").Value = -10911.0625
$ws.Range("H126").Value = 4289.6055
$ws.Range("I126").Value = 3816.0588
$ws.Range("J126").Value = 4672.952
$ws.Range("K126").Value = 11448.1764
$ws.Range("L126").Value = 14018.856
$ws.Range("M126").Value = -8978.1764
$ws.Range("N126").Value = -18958.856

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 853.25
$ws.Range("I113").Value = 306.66666
$ws.Range("J113").Value = 1556
$ws.Range("K113").Value = 919.9999799999999
$ws.Range("L113").Value = 4668
$ws.Range("M113").Value = 1250.00002
$ws.Range("N113").Value = -9008
$ws.Range("H126").Value = 3304.3242
$ws.Range("I126").Value = 4125.654
$ws.Range("J126").Value = 1363
$ws.Range("K126").Value = 12376.962
$ws.Range("L126").Value = 4089
$ws.Range("M126").Value = -9906.962000000001
$ws.Range("N126").Value = -9029
$ws.Range("H132").Value = 4387.3833
$ws.Range("I132").Value = 4751.22
$ws.Range("J132").Value = 2568.2
$ws.Range("K132").Value = 14253.66
$ws.Range("L132").Value = 7704.599999999999
$ws.Range("M132").Value = -11723.66
$ws.Range("N132").Value = -12764.6

Write-Host "All edits applied"